$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" --------------------------------------------------
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Publication date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Remove the duplicated "Contact / No display for ContactDetail" row
# (row 11 duplicates row 10); deleting it shifts everything below up by one
$ws.Rows.Item(11).Delete()

# Publisher now has a value
$ws.Range("B9").Value = "Alvearie Team"

# Row that used to be "Contact / No display for ContactDetail" becomes
# "Jurisdiction / United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# --- Sheet "Elements" ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("Elements")

# Top-level Extension row: Short/Definition now reflect the profile's own
# title/description instead of the generic Extension text
$ws2.Range("K2").Value = "Offset Begin"
$ws2.Range("L2").Value = "Offset location of the first character for the span of covered text in relation to the overall reference where this span of text appears"
